$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-21 down to 10-22
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new weekly record
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Femacal de La Calera"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44557
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 100112044
$ws.Range("G9").Value = "Perejil"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 104
$ws.Range("K9").Value = 2000
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2260
$ws.Range("N9").Value = "$/docena de atados (3 kilos)"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 753
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = "Hortaliza"
